$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: sr, workspace_id, repo_name_to_import, gitlab_target_namespace
$workspaceId = "anilgoudasb06"
$namespace   = "repo-migration"
$repos = @("almatasks", "app-n-pak", "casa-build-utils", "casaplotserver", "casashell")

# Write column B (workspace_id) first for every row, so the shared string
# "anilgoudasb06" is registered before the repo names.
for ($i = 0; $i -lt $repos.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $workspaceId
}

# Then write column C (repo_name_to_import) for every row, in order.
for ($i = 0; $i -lt $repos.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 3).Value = $repos[$i]
}

# Then write column D (gitlab_target_namespace) for every row.
for ($i = 0; $i -lt $repos.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 4).Value = $namespace
}

# Carry the existing D-column formatting down onto the newly added rows
# (D4:D6) so the whole column keeps a consistent style.
$ws.Range("D2").Copy()
$ws.Range("D4:D6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C4").Select() | Out-Null
